$d = $word.ActiveDocument
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n)
$insertPoint = $lastPara.Range.Duplicate
$insertPoint.Collapse(1)
$insertPoint.InsertParagraphBefore()

$slot = $d.Paragraphs.Item($n)
$slotRange = $slot.Range.Duplicate
$slotRange.Collapse(0)

$xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Nota: Ao que parece ele coloca o nome que a gente escrever no plural, como se fosse um conjunto daquelas coisas...</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">O </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Django</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> faz isso automaticamente, para ajustar vamos no arquivo models.py e vamos adicionar uma classe</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="1F1F1F"/><w:spacing w:after="0" w:line="285" w:lineRule="atLeast"/><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="CCCCCC"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="CCCCCC"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t xml:space="preserve">    </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="569CD6"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>class</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="CCCCCC"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="4EC9B0"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>Meta</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="CCCCCC"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>:</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="1F1F1F"/><w:spacing w:after="0" w:line="285" w:lineRule="atLeast"/><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="CCCCCC"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="CCCCCC"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t xml:space="preserve">        </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="9CDCFE"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>verbose</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="9CDCFE"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>_name_plural</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="CCCCCC"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="D4D4D4"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>=</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="CCCCCC"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="CE9178"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>'Categorias'</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="CCCCCC"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="6A9955"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t xml:space="preserve">#Aqui eu altero a forma como o </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="6A9955"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>Django</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="6A9955"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t xml:space="preserve"> vai transformar em plural as palavras que eu escrever.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Ainda no arquivo models.py vamos criar a classe Item, para representar os itens que o usuário quiser adicionar ao site, dentro de cada uma das categorias (são os produtos em si).</w:t></w:r></w:p><w:p><w:r><w:t>Depois de criarmos o campo Item, temos que ir no arquivo admin.py e registrar ele também, assim como foi com o campo Categoria.</w:t></w:r></w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$slotRange.InsertXML($xml)

# Now clear the text of the original last paragraph (now pushed further down), leaving only its bookmark.
$n2 = $d.Paragraphs.Count
$bmPara = $d.Paragraphs.Item($n2)
$full = $bmPara.Range.Duplicate
$txt = $full.Text
$len = $txt.Length
if ($len -gt 0) {
    $textOnly = $d.Range($full.Start, $full.Start + $len - 1)
    $textOnly.Text = ""
}

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
